# Update workbook with refreshed market price data (get available usd)
$wb = $excel.ActiveWorkbook

# --- Sheet 1: safety_orders -----------------------------------------------
$ws1 = $wb.Worksheets.Item("safety_orders")

$ws1.Range("D2").Value = 304.71762034176
$ws1.Range("E2").Value = 316.80026494488
$ws1.Range("F2").Value = 319.9682675943288
$ws1.Range("G2").Value = 4.766299910684991

$ws1.Range("D3").Value = 278.230297726265
$ws1.Range("E3").Value = 297.5152813355725
$ws1.Range("F3").Value = 300.4904341489282
$ws1.Range("G3").Value = 7.407935126357701

$ws1.Range("D4").Value = 236.910074459854
$ws1.Range("E4").Value = 267.2126778977132
$ws1.Range("F4").Value = 269.8848046766904
$ws1.Range("G4").Value = 12.21807587735019

$ws1.Range("D5").Value = 172.45052615462
$ws1.Range("E5").Value = 219.8316020261666
$ws1.Range("F5").Value = 222.0299180464283
$ws1.Range("G5").Value = 22.33005007975584

# --- Sheet 2: open_buy_orders ----------------------------------------------
$ws2 = $wb.Worksheets.Item("open_buy_orders")

$ws2.Range("A2").Value = "O736YO-27F72-JGBIPY"
$ws2.Range("B2").Value = 332.17

# Row 3 (OOA6WP-64JU5-6KOYXS / 321.79) no longer exists - delete it
$ws2.Rows.Item(3).Delete()

# --- Sheet 3: open_sell_orders ----------------------------------------------
$ws3 = $wb.Worksheets.Item("open_sell_orders")

$ws3.Range("A2").Value = "OGI3BN-DA3HY-C25WXB"

# Row 3 (OJG3BE-436LG-CLB5NB) no longer exists - delete it
$ws3.Rows.Item(3).Delete()
